$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper fragments used to build minimal WordprocessingML packages that can
# be fed to Range.InsertXML (which REPLACES the target range's contents).
# ---------------------------------------------------------------------------
$pkgHead = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>'
$pkgTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1) Empty paragraph that used to host the _GoBack bookmark. -----------
# Drop the bookmark from there, and add w:hint="cs", w:rtl and
# w:lang w:bidi="ar-EG" to the (empty) paragraph mark run properties.
$goBackPara = $null
foreach ($p in $d.Paragraphs) {
    $pxml = $p.Range.WordOpenXML
    if ($pxml -like "*_GoBack*") {
        $goBackPara = $p
        break
    }
}

$p3Xml = '<w:p w:rsidR="00AE06F2" w:rsidRDefault="00AE06F2" w:rsidP="00F12605">' +
    '<w:pPr>' +
    '<w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>' +
    '<w:outlineLvl w:val="0"/>' +
    '<w:rPr>' +
    '<w:rFonts w:asciiTheme="majorBidi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi" w:hint="cs"/>' +
    '<w:b/><w:bCs/><w:kern w:val="36"/><w:sz w:val="48"/><w:szCs w:val="48"/>' +
    '<w:rtl/><w:lang w:bidi="ar-EG"/>' +
    '</w:rPr>' +
    '</w:pPr>' +
    '</w:p>'

$p3Pkg = $pkgHead + $p3Xml + $pkgTail
$goBackPara.Range.InsertXML($p3Pkg)

# --- 2) Split the "Identify critical defects..." sentence into 4 runs and
# re-insert the _GoBack bookmark so that it now wraps "17" (the date changes
# from May 24 to May 17). -----------------------------------------------
$oldSentence = "Identify critical defects before the release on Saturday, May 24, 2025."
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ($oldSentence + "*")) {
        $targetPara = $p
        break
    }
}

$p19Xml = '<w:p w:rsidR="00F12605" w:rsidRPr="00DF7120" w:rsidRDefault="00F12605" w:rsidP="00F12605">' +
    '<w:pPr>' +
    '<w:pStyle w:val="NormalWeb"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="00DF7120">' +
    '<w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '<w:t xml:space="preserve">Identify critical defects before </w:t>' +
    '</w:r>' +
    '<w:r>' +
    '<w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '<w:t xml:space="preserve">the release on Saturday, May </w:t>' +
    '</w:r>' +
    '<w:r>' +
    '<w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/></w:rPr>' +
    '<w:t>17</w:t>' +
    '</w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r>' +
    '<w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
    '<w:t>, 2025.</w:t>' +
    '</w:r>' +
    '</w:p>'

$p19Pkg = $pkgHead + $p19Xml + $pkgTail
$targetPara.Range.InsertXML($p19Pkg)
